$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 120, pushing the existing rows 120:140 down to 121:141.
$ws.Rows.Item(120).Insert()

# Populate the new row 120 with a fresh weekly price entry for
# "Agrícola del Norte S.A. de Arica" / "Cebollín baby" (same as the prior
# top row of this block, but dated one period later).
$ws.Range("A120").Value = 1
$ws.Range("B120").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C120").Value = "Arica y Parinacota"
$ws.Range("D120").Value = 45173
$ws.Range("D120").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E120").Value = 15
$ws.Range("F120").Value = 100112038
$ws.Range("G120").Value = "Cebollín baby"
$ws.Range("H120").Value = "Sin especificar"
$ws.Range("I120").Value = "Primera"
$ws.Range("J120").Value = 300
$ws.Range("K120").Value = 1900
$ws.Range("L120").Value = 2000
$ws.Range("M120").Value = 1950
$ws.Range("N120").Value = "$/paquete 1,5 a 2 kilos"
$ws.Range("O120").Value = "Región de Arica y Parinacota"
$ws.Range("P120").Value = 975
$ws.Range("Q120").Value = 2
$ws.Range("R120").Value = "Hortaliza"
